$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows 83-88: "QUESTION" list continuation (col A = serial no, col B = question)
$ws.Cells.Item(83, 1).Value = 6
$ws.Cells.Item(83, 2).Value = "Check if array is Sorted"

$ws.Cells.Item(84, 1).Value = 7
$ws.Cells.Item(84, 2).Value = "Binary Search"

$ws.Cells.Item(85, 1).Value = 8
$ws.Cells.Item(85, 2).Value = "Print All Subsequences"

$ws.Cells.Item(86, 1).Value = 9
$ws.Cells.Item(86, 2).Value = "Minimum no of elements whose sum=target"

$ws.Cells.Item(87, 1).Value = 10
$ws.Cells.Item(87, 2).Value = "Cut into Segments"

$ws.Cells.Item(88, 1).Value = 11
$ws.Cells.Item(88, 2).Value = "Maximum sum by using  non adjacent elements"

# Match the formatting used by the rest of the "S.No" column (A) above
$ws.Range("A82").Copy()
$ws.Range("A83:A88").PasteSpecial(-4122)

# Last new row's question cell gets vertical-center alignment
$ws.Cells.Item(88, 2).VerticalAlignment = -4108

# Match row heights used throughout the sheet
$ws.Range("A83:B88").RowHeight = 15.6

# Move the view/selection down to the newly added rows
$ws.Activate()
$ws.Range("C87").Select()
